$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031467113740999
$ws.Range("D2").Value = 1.035531834160866
$ws.Range("E2").Value = 1.051332613646209
$ws.Range("F2").Value = 1.056530000849288
$ws.Range("I2").Value = 1.035321186475805
$ws.Range("J2").Value = 1.036602680376708
$ws.Range("K2").Value = 1.038328170421941
$ws.Range("L2").Value = 1.054084408225717
$ws.Range("M2").Value = 1.059267459768043
$ws.Range("N2").Value = 1.038074775082832

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03230456995521
$ws.Range("D3").Value = 1.036156212083558
$ws.Range("E3").Value = 1.052555716260776
$ws.Range("F3").Value = 1.057813192493796
$ws.Range("I3").Value = 1.035493991641943
$ws.Range("J3").Value = 1.037082594133174
$ws.Range("K3").Value = 1.038762454573666
$ws.Range("L3").Value = 1.055119080381258
$ws.Range("M3").Value = 1.060363113860959
$ws.Range("N3").Value = 1.038555370371879

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03284664574786
$ws.Range("D4").Value = 1.036560334472244
$ws.Range("E4").Value = 1.053348056503081
$ws.Range("F4").Value = 1.058644448253293
$ws.Range("I4").Value = 1.035604701691443
$ws.Range("J4").Value = 1.037392660709816
$ws.Range("K4").Value = 1.039042882967115
$ws.Range("L4").Value = 1.055788927128594
$ws.Range("M4").Value = 1.061072464738002
$ws.Range("N4").Value = 1.038865877278625

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033074577807439
$ws.Range("D5").Value = 1.036730252095992
$ws.Range("E5").Value = 1.053681375292724
$ws.Range("F5").Value = 1.05899413523252
$ws.Range("I5").Value = 1.035650978843097
$ws.Range("J5").Value = 1.037522899481397
$ws.Range("K5").Value = 1.039160634759524
$ws.Range("L5").Value = 1.05607061391045
$ws.Range("M5").Value = 1.061370769700054
$ws.Range("N5").Value = 1.0389963010042

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033112851117708
$ws.Range("D6").Value = 1.036758783415346
$ws.Range("E6").Value = 1.053737353848258
$ws.Range("F6").Value = 1.05905286255421
$ws.Range("I6").Value = 1.035658733406488
$ws.Range("J6").Value = 1.037544760482822
$ws.Range("K6").Value = 1.039180397551244
$ws.Range("L6").Value = 1.056117915260355
$ws.Range("M6").Value = 1.061420861937784
$ws.Range("N6").Value = 1.039018193050755

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032849691217979
$ws.Range("D7").Value = 1.036562604824954
$ws.Range("E7").Value = 1.05335250946366
$ws.Range("F7").Value = 1.058649119894582
$ws.Range("I7").Value = 1.035605321092246
$ws.Range("J7").Value = 1.037394401412483
$ws.Range("K7").Value = 1.0390444569247
$ws.Range("L7").Value = 1.055792690713323
$ws.Range("M7").Value = 1.061076450333927
$ws.Range("N7").Value = 1.03886762045329

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031750096928357
$ws.Range("D8").Value = 1.035742822427427
$ws.Range("E8").Value = 1.051745778763438
$ws.Range("F8").Value = 1.056963466311568
$ws.Range("I8").Value = 1.03537981560708
$ws.Range("J8").Value = 1.0367649662745
$ws.Range("K8").Value = 1.038475058879967
$ws.Range("L8").Value = 1.05443400957258
$ws.Range("M8").Value = 1.059637660942777
$ws.Range("N8").Value = 1.038237291445214

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029813938838468
$ws.Range("D9").Value = 1.034299147006819
$ws.Range("E9").Value = 1.048921454020306
$ws.Range("F9").Value = 1.054000329354389
$ws.Range("I9").Value = 1.034973992667583
$ws.Range("J9").Value = 1.035652257410964
$ws.Range("K9").Value = 1.037467283186701
$ws.Range("L9").Value = 1.052042448600219
$ws.Range("M9").Value = 1.057105284571296
$ws.Range("N9").Value = 1.037123002407462

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028524214065362
$ws.Range("D10").Value = 1.033337369678381
$ws.Range("E10").Value = 1.047043182227397
$ws.Range("F10").Value = 1.052029691262703
$ws.Range("I10").Value = 1.034697787801856
$ws.Range("J10").Value = 1.034908100073474
$ws.Range("K10").Value = 1.036792508569441
$ws.Range("L10").Value = 1.050449783318762
$ws.Range("M10").Value = 1.055418978568751
$ws.Range("N10").Value = 1.036377788281197

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027966009765311
$ws.Range("D11").Value = 1.032921085198608
$ws.Range("E11").Value = 1.046230950705684
$ws.Range("F11").Value = 1.051177505696923
$ws.Range("I11").Value = 1.03457685209397
$ws.Range("J11").Value = 1.034585323132806
$ws.Range("K11").Value = 1.036499640458303
$ws.Range("L11").Value = 1.049760536241487
$ws.Range("M11").Value = 1.05468924084015
$ws.Range("N11").Value = 1.03605455296025

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027758707095373
$ws.Range("D12").Value = 1.03276648564201
$ws.Range("E12").Value = 1.045929411703011
$ws.Range("F12").Value = 1.050861132166141
$ws.Range("I12").Value = 1.034531730714429
$ws.Range("J12").Value = 1.034465347183891
$ws.Range("K12").Value = 1.036390753897634
$ws.Range("L12").Value = 1.049504576401472
$ws.Range("M12").Value = 1.054418249503041
$ws.Range("N12").Value = 1.035934406631725

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027803172429956
$ws.Range("D13").Value = 1.0327996465279
$ws.Range("E13").Value = 1.045994085623704
$ws.Range("F13").Value = 1.050928987863728
$ws.Range("I13").Value = 1.034541418468533
$ws.Range("J13").Value = 1.034491086141734
$ws.Range("K13").Value = 1.036414115051221
$ws.Range("L13").Value = 1.049559478083074
$ws.Range("M13").Value = 1.054476375104659
$ws.Range("N13").Value = 1.035960182141841

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027948873255563
$ws.Range("D14").Value = 1.032908305392178
$ws.Range("E14").Value = 1.046206022141606
$ws.Range("F14").Value = 1.051151350780012
$ws.Range("I14").Value = 1.0345731264338
$ws.Range("J14").Value = 1.034575407556715
$ws.Range("K14").Value = 1.036490641943602
$ws.Range("L14").Value = 1.049739377359916
$ws.Range("M14").Value = 1.054666839275907
$ws.Range("N14").Value = 1.036044623302903

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028038649573711
$ws.Range("D15").Value = 1.032975257325051
$ws.Range("E15").Value = 1.046336624384381
$ws.Range("F15").Value = 1.051288377890145
$ws.Range("I15").Value = 1.034592636205749
$ws.Range("J15").Value = 1.034627349889095
$ws.Range("K15").Value = 1.036537779155217
$ws.Range("L15").Value = 1.049850226793186
$ws.Range("M15").Value = 1.054784199221905
$ws.Range("N15").Value = 1.036096639399353

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028561265700316
$ws.Range("D16").Value = 1.033365000829551
$ws.Range("E16").Value = 1.047097109854272
$ws.Range("F16").Value = 1.052086271375103
$ws.Range("I16").Value = 1.034705785769767
$ws.Range("J16").Value = 1.03492951014742
$ws.Range("K16").Value = 1.036811930871012
$ws.Range("L16").Value = 1.05049553445561
$ws.Range("M16").Value = 1.055467418130184
$ws.Range("N16").Value = 1.036399228759904

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028889157936669
$ws.Range("D17").Value = 1.033609523437014
$ws.Range("E17").Value = 1.047574428610048
$ws.Range("F17").Value = 1.052587066211389
$ws.Range("I17").Value = 1.034776403799626
$ws.Range("J17").Value = 1.035118900025763
$ws.Range("K17").Value = 1.036983715833988
$ws.Range("L17").Value = 1.050900421937532
$ws.Range("M17").Value = 1.055896101477183
$ws.Range("N17").Value = 1.036588887593597

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029080436474686
$ws.Range("D18").Value = 1.033752165794068
$ws.Range("E18").Value = 1.04785294413064
$ws.Range("F18").Value = 1.052879278905315
$ws.Range("I18").Value = 1.034817465054948
$ws.Range("J18").Value = 1.035229314640339
$ws.Range("K18").Value = 1.037083848830065
$ws.Range("L18").Value = 1.051136623410684
$ws.Range("M18").Value = 1.056146188113453
$ws.Range("N18").Value = 1.03669945900959

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029145661647436
$ws.Range("D19").Value = 1.033800805904154
$ws.Range("E19").Value = 1.047947928328574
$ws.Range("F19").Value = 1.052978934159333
$ws.Range("I19").Value = 1.034831443980475
$ws.Range("J19").Value = 1.03526695410923
$ws.Range("K19").Value = 1.037117980329855
$ws.Range("L19").Value = 1.051217168392712
$ws.Range("M19").Value = 1.056231468497758
$ws.Range("N19").Value = 1.036737151930844

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028853975657276
$ws.Range("D20").Value = 1.033583286763896
$ws.Range("E20").Value = 1.047523206120966
$ws.Range("F20").Value = 1.052533324550138
$ws.Range("I20").Value = 1.034768840498186
$ws.Range("J20").Value = 1.03509858579191
$ws.Range("K20").Value = 1.036965291781685
$ws.Range("L20").Value = 1.050856977475481
$ws.Range("M20").Value = 1.055850103381133
$ws.Range("N20").Value = 1.036568544511203

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027905966902529
$ws.Range("D21").Value = 1.032876307308394
$ws.Range("E21").Value = 1.046143607713617
$ws.Range("F21").Value = 1.051085865866916
$ws.Range("I21").Value = 1.034563794758963
$ws.Range("J21").Value = 1.034550579265791
$ws.Range("K21").Value = 1.036468109498882
$ws.Range("L21").Value = 1.049686399935427
$ws.Range("M21").Value = 1.054610750517282
$ws.Range("N21").Value = 1.036019759752958

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027310143793267
$ws.Range("D22").Value = 1.032431958305411
$ws.Range("E22").Value = 1.045277124270229
$ws.Range("F22").Value = 1.050176751745346
$ws.Range("I22").Value = 1.034433714699707
$ws.Range("J22").Value = 1.034205550448199
$ws.Range("K22").Value = 1.036154919757433
$ws.Range("L22").Value = 1.048950741901343
$ws.Range("M22").Value = 1.053831900110513
$ws.Range("N22").Value = 1.03567424095487

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027625979005321
$ws.Range("D23").Value = 1.03266750080863
$ws.Range("E23").Value = 1.045736376270693
$ws.Range("F23").Value = 1.050658599685171
$ws.Range("I23").Value = 1.034502782381679
$ws.Range("J23").Value = 1.034388501579478
$ws.Range("K23").Value = 1.036321003367386
$ws.Range("L23").Value = 1.049340697148122
$ws.Range("M23").Value = 1.054244747748251
$ws.Range("N23").Value = 1.035857451897741

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02886987293915
$ws.Range("D24").Value = 1.033595141936413
$ws.Range("E24").Value = 1.047546351038572
$ws.Range("F24").Value = 1.052557607760925
$ws.Range("I24").Value = 1.03477225842731
$ws.Range("J24").Value = 1.035107765084494
$ws.Range("K24").Value = 1.036973617022681
$ws.Range("L24").Value = 1.050876608041788
$ws.Range("M24").Value = 1.055870887807924
$ws.Range("N24").Value = 1.036577736839434

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030314301961396
$ws.Range("D25").Value = 1.034672258742623
$ws.Range("E25").Value = 1.049650792895054
$ws.Range("F25").Value = 1.054765523666331
$ws.Range("I25").Value = 1.035079906584244
$ws.Range("J25").Value = 1.035940337043507
$ws.Range("K25").Value = 1.037728336419294
$ws.Range("L25").Value = 1.052660419970556
$ws.Range("M25").Value = 1.057759618747437
$ws.Range("N25").Value = 1.037411491146127
